$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1749
$ws1.Range("F5").Value = 440
$ws1.Range("F9").Value = 295
$ws1.Range("F10").Value = 1678
$ws1.Range("F11").Value = 334
$ws1.Range("F12").Value = 1387
$ws1.Range("F13").Value = 786
$ws1.Range("F14").Value = 318
$ws1.Range("F15").Value = 654
$ws1.Range("F16").Value = 12617
$ws1.Range("F17").Value = 12645
$ws1.Range("F18").Value = 933
$ws1.Range("F19").Value = 732
$ws1.Range("F21").Value = 297
$ws1.Range("F22").Value = 45
$ws1.Range("F23").Value = 491
$ws1.Range("F27").Value = 228
$ws1.Range("F28").Value = 661

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1749
$ws4.Range("F7").Value = 440
$ws4.Range("F14").Value = 295
$ws4.Range("F15").Value = 1678
$ws4.Range("F16").Value = 334
$ws4.Range("F17").Value = 1387
$ws4.Range("F18").Value = 786
$ws4.Range("F19").Value = 318
$ws4.Range("F21").Value = 654
$ws4.Range("F22").Value = 12617
$ws4.Range("F23").Value = 12645
$ws4.Range("F24").Value = 933
$ws4.Range("F25").Value = 732
$ws4.Range("F27").Value = 297
$ws4.Range("F28").Value = 45
$ws4.Range("F29").Value = 491
$ws4.Range("F37").Value = 228
$ws4.Range("F38").Value = 661
